# Updated cryptos list - price (D) and volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.439.90"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.02"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.88"
$ws.Range("E5").Value = "  -5.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5203"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3275"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.63"
$ws.Range("E10").Value = "  -7.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7720"
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07695"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.823.39"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.20"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.026"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.93"
$ws.Range("E17").Value = "  -4.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007953"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.486.21"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.071.45"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.578"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.483"
$ws.Range("E23").Value = "  -5.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.970"
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.26"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.205"
$ws.Range("E26").Value = "  -7.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.655"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.56"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.194"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.131"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08721"
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04795"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.129"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7097"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.836"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.093"
$ws.Range("E37").Value = "  -6.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.240"
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01763"
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4839"
$ws.Range("E40").Value = "  -5.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.83"
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8929"
$ws.Range("E42").Value = "  -6.55%  "
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.724"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4163"
$ws.Range("E46").Value = "  -6.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05869"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.004"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.97"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1217"
$ws.Range("E50").Value = "  -9.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8878"
$ws.Range("E51").Value = "  +0.48%  "
